$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2245.652
$ws.Range("I132").Value = 2120.4546
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 6361.3638
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -3831.3638
$ws.Range("N132").Value = -20060
$ws.Range("H135").Value = 1387.6072
$ws.Range("I135").Value = 988.9474
$ws.Range("K135").Value = 8900.526600000001
$ws.Range("M135").Value = -6365.526600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1422.875
$ws.Range("I2").Value = 1235.1666
$ws.Range("J2").Value = 1986
$ws.Range("K2").Value = 1235.1666
$ws.Range("L2").Value = 1986
$ws.Range("M2").Value = -1122.1666
$ws.Range("N2").Value = -2212
$ws.Range("H32").Value = 9770.759
$ws.Range("I32").Value = 4871.381
$ws.Range("J32").Value = 26918.584
$ws.Range("K32").Value = 4871.381
$ws.Range("L32").Value = 26918.584
$ws.Range("M32").Value = -4584.381
$ws.Range("N32").Value = -27492.584
$ws.Range("H61").Value = 3879.5293
$ws.Range("I61").Value = 2870.5
$ws.Range("K61").Value = 2870.5
$ws.Range("M61").Value = -2658.5
$ws.Range("H74").Value = 9479.280000000001
$ws.Range("I74").Value = 1689.4736
$ws.Range("K74").Value = 1689.4736
$ws.Range("M74").Value = -815.4736
$ws.Range("H77").Value = 9479.280000000001
$ws.Range("I77").Value = 1689.4736
$ws.Range("K77").Value = 8447.368
$ws.Range("M77").Value = -4079.368
$ws.Range("H97").Value = 45095.26
$ws.Range("I97").Value = 1184.7368
$ws.Range("J97").Value = 253670.25
$ws.Range("K97").Value = 1184.7368
$ws.Range("L97").Value = 253670.25
$ws.Range("M97").Value = -688.7367999999999
$ws.Range("N97").Value = -254662.25
$ws.Range("H110").Value = 6884.1577
$ws.Range("I110").Value = 7364.647
$ws.Range("K110").Value = 7364.647
$ws.Range("M110").Value = -5319.647
$ws.Range("H116").Value = 1422.875
$ws.Range("I116").Value = 1235.1666
$ws.Range("J116").Value = 1986
$ws.Range("K116").Value = 1235.1666
$ws.Range("L116").Value = 1986
$ws.Range("M116").Value = 1058.8334
$ws.Range("N116").Value = -6574
$ws.Range("H132").Value = 1600.1842
$ws.Range("I132").Value = 1355.75
$ws.Range("K132").Value = 4067.25
$ws.Range("M132").Value = -1537.25
$ws.Range("H136").Value = 3879.5293
$ws.Range("I136").Value = 2870.5
$ws.Range("K136").Value = 8611.5
$ws.Range("M136").Value = -6061.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1422.875
$ws.Range("I3").Value = 1235.1666
$ws.Range("J3").Value = 1986
$ws.Range("K3").Value = 1235.1666
$ws.Range("L3").Value = 1986
$ws.Range("M3").Value = -1121.1666
$ws.Range("N3").Value = -2214
$ws.Range("H107").Value = 1111.5238
$ws.Range("I107").Value = 1089
$ws.Range("J107").Value = 1246.6666
$ws.Range("K107").Value = 1089
$ws.Range("L107").Value = 1246.6666
$ws.Range("M107").Value = 831
$ws.Range("N107").Value = -5086.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 2492.5557
$ws.Range("I94").Value = 1408.25
$ws.Range("K94").Value = 1408.25
$ws.Range("M94").Value = -957.25
$ws.Range("H105").Value = 1504.8
$ws.Range("I105").Value = 1504.8
$ws.Range("K105").Value = 1504.8
$ws.Range("M105").Value = 242.2
$ws.Range("H132").Value = 2173.5227
$ws.Range("I132").Value = 1969.05
$ws.Range("K132").Value = 5907.15
$ws.Range("M132").Value = -3377.15

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 266.6875
$ws.Range("J12").Value = 266.2857
$ws.Range("L12").Value = 798.8571000000001
$ws.Range("N12").Value = -1144.8571
$ws.Range("H32").Value = 7508.3335
$ws.Range("J32").Value = 9429.333000000001
$ws.Range("L32").Value = 28287.999
$ws.Range("N32").Value = -28853.999
$ws.Range("H46").Value = 333583.34
$ws.Range("I46").Value = 500150
$ws.Range("K46").Value = 1500450
$ws.Range("M46").Value = -1500359
$ws.Range("H75").Value = 2150
$ws.Range("I75").Value = 2725
$ws.Range("J75").Value = 1000
$ws.Range("K75").Value = 8175
$ws.Range("L75").Value = 3000
$ws.Range("M75").Value = -7177
$ws.Range("N75").Value = -4996
$ws.Range("H78").Value = 2150
$ws.Range("I78").Value = 2725
$ws.Range("J78").Value = 1000
$ws.Range("K78").Value = 24525
$ws.Range("L78").Value = 9000
$ws.Range("M78").Value = -19533
$ws.Range("N78").Value = -18984
$ws.Range("H114").Value = 2260
$ws.Range("J114").Value = 3224.5
$ws.Range("L114").Value = 9673.5
$ws.Range("N114").Value = -16181.5
$ws.Range("H121").Value = 50514.5
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H134").Value = 7012.375
$ws.Range("I134").Value = 4683.1665
$ws.Range("K134").Value = 14049.4995
$ws.Range("M134").Value = -8979.499500000002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 15062.75
$ws.Range("J92").Value = 20125.5
$ws.Range("L92").Value = 20125.5
$ws.Range("N92").Value = -23869.5
$ws.Range("H98").Value = 10321
$ws.Range("J98").Value = 10321
$ws.Range("L98").Value = 10321
$ws.Range("N98").Value = -16311
$ws.Range("H107").Value = 1169.2
$ws.Range("I107").Value = 1111.75
$ws.Range("J107").Value = 1399
$ws.Range("K107").Value = 1111.75
$ws.Range("L107").Value = 1399
$ws.Range("M107").Value = 808.25
$ws.Range("N107").Value = -5239
$ws.Range("H113").Value = 2850.7778
$ws.Range("I113").Value = 3101
$ws.Range("J113").Value = 1975
$ws.Range("K113").Value = 3101
$ws.Range("L113").Value = 1975
$ws.Range("M113").Value = -931
$ws.Range("N113").Value = -6315
$ws.Range("H132").Value = 2479.8462
$ws.Range("I132").Value = 2152.697
$ws.Range("J132").Value = 4279.1665
$ws.Range("K132").Value = 6458.091
$ws.Range("L132").Value = 12837.4995
$ws.Range("M132").Value = -3928.091
$ws.Range("N132").Value = -17897.4995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5687.75
$ws.Range("J16").Value = 7616.3335
$ws.Range("L16").Value = 7616.3335
$ws.Range("N16").Value = -7956.3335
$ws.Range("H46").Value = 2550
$ws.Range("I46").Value = 1680
$ws.Range("K46").Value = 1680
$ws.Range("M46").Value = -1492
$ws.Range("H61").Value = 4981.75
$ws.Range("I61").Value = 4798.273
$ws.Range("K61").Value = 4798.273
$ws.Range("M61").Value = -4596.273
$ws.Range("H113").Value = 4981.75
$ws.Range("I113").Value = 4798.273
$ws.Range("K113").Value = 4798.273
$ws.Range("M113").Value = -2628.273
$ws.Range("H132").Value = 5166.846
$ws.Range("I132").Value = 4417
$ws.Range("J132").Value = 7666.3335
$ws.Range("K132").Value = 13251
$ws.Range("L132").Value = 22999.0005
$ws.Range("M132").Value = -10721
$ws.Range("N132").Value = -28059.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H131").Value = 86571.28999999999
$ws.Range("J131").Value = 86571.28999999999
$ws.Range("L131").Value = 86571.28999999999
$ws.Range("N131").Value = -96651.28999999999
$ws.Range("H132").Value = 2034.25
$ws.Range("I132").Value = 1788.0731
$ws.Range("K132").Value = 5364.219300000001
$ws.Range("M132").Value = -2834.219300000001
$ws.Range("H136").Value = 1657.25
$ws.Range("I136").Value = 1449.4138
$ws.Range("J136").Value = 3666.3333
$ws.Range("K136").Value = 4348.2414
$ws.Range("L136").Value = 10998.9999
$ws.Range("M136").Value = -1798.2414
$ws.Range("N136").Value = -16098.9999
